$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data revision: flag column F (edu8) from 0 -> 1 for the rows whose
# underlying source data changed.
$rows = 36,39,43,48,51,52,53,54,61,67,68,84,97,110,112,114,122,124,147,150,158,162,166,173,175,181,185,195,203,204,205,221,223
foreach ($r in $rows) {
    $ws.Range("F$r").Value = 1
}

# Reflect the author's final scroll position / selection in the sheet view.
$ws.Range("I179").Select()
